# Auto-generated Excel COM-interop script applying the cryptos.xlsx price/volume update
# described by the commit "Updated cryptos list on Tue May  9 15:53:16 UTC 2023 with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain TEXT (matches the workbook's existing inline-string cells)
# without leaving a stray custom number-format style behind on the cell.
function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.558.05"
Set-TextValue $ws.Range("E2") "  -1.50%  "
Set-TextValue $ws.Range("D3") "1.847.83"
Set-TextValue $ws.Range("E3") "  -1.26%  "
Set-TextValue $ws.Range("D4") "1.004"
Set-TextValue $ws.Range("E4") "  +0.12%  "
Set-TextValue $ws.Range("D5") "314.10"
Set-TextValue $ws.Range("E5") "  -1.42%  "
Set-TextValue $ws.Range("E6") "  +0.05%  "
Set-TextValue $ws.Range("D7") "0.4243"
Set-TextValue $ws.Range("E7") "  -2.54%  "
Set-TextValue $ws.Range("D8") "0.3631"
Set-TextValue $ws.Range("E8") "  -3.36%  "
Set-TextValue $ws.Range("E9") "  +0.39%  "
Set-TextValue $ws.Range("D10") "0.07289"
Set-TextValue $ws.Range("E10") "  -2.50%  "
Set-TextValue $ws.Range("D11") "0.8704"
Set-TextValue $ws.Range("E11") "  -7.16%  "
Set-TextValue $ws.Range("D12") "20.67"
Set-TextValue $ws.Range("E12") "  -2.81%  "
Set-TextValue $ws.Range("D13") "1.901.84"
Set-TextValue $ws.Range("E13") "  -1.40%  "
Set-TextValue $ws.Range("D14") "5.337"
Set-TextValue $ws.Range("E14") "  -1.83%  "
Set-TextValue $ws.Range("D15") "6.504"
Set-TextValue $ws.Range("E15") "  -3.60%  "
Set-TextValue $ws.Range("D16") "0.06923"
Set-TextValue $ws.Range("E16") "  +0.79%  "
Set-TextValue $ws.Range("D17") "1.006"
Set-TextValue $ws.Range("E17") "  +0.28%  "
Set-TextValue $ws.Range("D18") "78.78"
Set-TextValue $ws.Range("E18") "  -3.25%  "
Set-TextValue $ws.Range("D19") "0.000008864"
Set-TextValue $ws.Range("E19") "  -2.10%  "
Set-TextValue $ws.Range("E20") "  +0.06%  "
Set-TextValue $ws.Range("D21") "15.38"
Set-TextValue $ws.Range("E21") "  -2.92%  "
Set-TextValue $ws.Range("D22") "27.588.67"
Set-TextValue $ws.Range("E22") "  -1.36%  "
Set-TextValue $ws.Range("D23") "5.000"
Set-TextValue $ws.Range("E23") "  -2.56%  "
Set-TextValue $ws.Range("D24") "10.60"
Set-TextValue $ws.Range("E24") "  -4.19%  "
Set-TextValue $ws.Range("D25") "2.125.21"
Set-TextValue $ws.Range("E25") "  +0.41%  "
Set-TextValue $ws.Range("D26") "1.983"
Set-TextValue $ws.Range("E26") "  -2.73%  "
Set-TextValue $ws.Range("D27") "153.54"
Set-TextValue $ws.Range("E27") "  +0.26%  "
Set-TextValue $ws.Range("D28") "18.96"
Set-TextValue $ws.Range("E28") "  +2.14%  "
Set-TextValue $ws.Range("D29") "120.92"
Set-TextValue $ws.Range("E29") "  +6.57%  "
Set-TextValue $ws.Range("D30") "5.259"
Set-TextValue $ws.Range("E30") "  -5.77%  "
Set-TextValue $ws.Range("D31") "1.902"
Set-TextValue $ws.Range("E31") "  +11.71%  "
Set-TextValue $ws.Range("D32") "0.08908"
Set-TextValue $ws.Range("E32") "  -1.41%  "
Set-TextValue $ws.Range("D33") "0.7587"
Set-TextValue $ws.Range("E33") "  -6.50%  "
Set-TextValue $ws.Range("D34") "4.560"
Set-TextValue $ws.Range("E34") "  -5.14%  "
Set-TextValue $ws.Range("E35") "  -1.08%  "
Set-TextValue $ws.Range("D36") "1.096"
Set-TextValue $ws.Range("E36") "  -7.09%  "
Set-TextValue $ws.Range("D37") "1.002"
Set-TextValue $ws.Range("E37") "  -0.05%  "
Set-TextValue $ws.Range("B38") "TrustWalletToken"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D38") "1.094"
Set-TextValue $ws.Range("E38") "  -2.46%  "
Set-TextValue $ws.Range("B39") "Hedera"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D39") "0.05359"
Set-TextValue $ws.Range("E39") "  -2.96%  "
Set-TextValue $ws.Range("D40") "0.01939"
Set-TextValue $ws.Range("E40") "  -2.29%  "
Set-TextValue $ws.Range("D41") "2.803"
Set-TextValue $ws.Range("E41") "  -5.95%  "
Set-TextValue $ws.Range("D42") "6.919"
Set-TextValue $ws.Range("E42") "  -1.07%  "
Set-TextValue $ws.Range("D43") "0.5100"
Set-TextValue $ws.Range("E43") "  -3.28%  "
Set-TextValue $ws.Range("D44") "0.1648"
Set-TextValue $ws.Range("E44") "  -3.09%  "
Set-TextValue $ws.Range("D45") "8.274"
Set-TextValue $ws.Range("E45") "  -5.87%  "
Set-TextValue $ws.Range("D46") "0.06559"
Set-TextValue $ws.Range("E46") "  -2.86%  "
Set-TextValue $ws.Range("B47") "EnergySwap"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D47") "10.42"
Set-TextValue $ws.Range("E47") "  -1.01%  "
Set-TextValue $ws.Range("B48") "Decentraland"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D48") "0.4745"
Set-TextValue $ws.Range("E48") "  -3.04%  "
Set-TextValue $ws.Range("D49") "104.32"
Set-TextValue $ws.Range("E49") "  -2.48%  "
Set-TextValue $ws.Range("D50") "1.002"
Set-TextValue $ws.Range("E50") "  +0.01%  "
Set-TextValue $ws.Range("D51") "1.621"
Set-TextValue $ws.Range("E51") "  -3.13%  "
